$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Append a new sentence after "..., a project management software."
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("a project management software.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter(" We also included some basic helper components, such as primitive types and exception handling from Lars's private code base.")

# Re-find across the run boundary so the new run's properties are
# addressable, then stamp its language explicitly (matches sibling runs).
$langFix = $d.Content
$langFix.Find.Execute("software. We also included", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$langFix.LanguageID = "en-US"

# ---------------------------------------------------------------------
# 2) British -> American spelling: "behaviour." -> "behavior." split
#    across three runs ("...beha" | "vio" | "r.") with no proofErr tags.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("their own behaviour.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Text = "their own behavior."

$rng2 = $d.Content
$rng2.Find.Execute("their own behavior.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$base = $rng2.Start

$vioRng = $d.Range($base + 14, $base + 17)
$vioRng.Bold = 1
$vioRng.Bold = 0

$rRng = $d.Range($base + 17, $base + 19)
$rRng.Bold = 1
$rRng.Bold = 0

# ---------------------------------------------------------------------
# 3) Append a new sentence after "... was added."
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("A task for Redmine was added.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter(" The render batch also needs to be postponed, since it cannot be tested without a scene.")

$langFix2 = $d.Content
$langFix2.Find.Execute("added. The render batch", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$langFix2.LanguageID = "en-US"

# ---------------------------------------------------------------------
# 4) The empty paragraph right after "...are shown below:" gains an
#    explicit language mark on its paragraph mark.
# ---------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("collaboratively, are shown below:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$emptyPara = $anchor.Paragraphs(1).Next()
$emptyPara.Range.LanguageID = "en-US"
